$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update citation in column A for rows 2 through 61
for ($r = 2; $r -le 61; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "Got2021, 34835480") {
        $cell.Value = "Miti2020, 32804970"
    }
}

# Update specific question text in column C
$ws.Range("C5").Value = "Does the paper report novel in vitro antiretroviral susceptibility data?"
$ws.Range("C55").Value = "Does the paper report novel IC values like IC50? IC90?"
$ws.Range("C56").Value = "Does the paper report novel IC50 fold change values?"

$wb.Save()
